$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 214.36363  # H6: 260.44446 -> 214.36363
$ws.Cells.Item(6, 9).Value = 135.8  # I6: 168 -> 135.8
$ws.Cells.Item(6, 11).Value = 407.4  # K6: 504 -> 407.4
$ws.Cells.Item(6, 13).Value = -295.4  # M6: -392 -> -295.4
$ws.Cells.Item(34, 8).Value = 19364.334  # H34: 14014.667 -> 19364.334
$ws.Cells.Item(34, 9).Value = 19022  # I34: 14014.667 -> 19022
$ws.Cells.Item(34, 10).Value = 20049  # J34: 0 -> 20049
$ws.Cells.Item(34, 11).Value = 19022  # K34: 14014.667 -> 19022
$ws.Cells.Item(34, 12).Value = 20049  # L34: 0 -> 20049
$ws.Cells.Item(34, 13).Value = -18819  # M34: -13811.667 -> -18819
$ws.Cells.Item(34, 14).Value = -20455  # N34: None -> -20455
$ws.Cells.Item(36, 8).Value = 19364.334  # H36: 14014.667 -> 19364.334
$ws.Cells.Item(36, 9).Value = 19022  # I36: 14014.667 -> 19022
$ws.Cells.Item(36, 10).Value = 20049  # J36: 0 -> 20049
$ws.Cells.Item(36, 11).Value = 19022  # K36: 14014.667 -> 19022
$ws.Cells.Item(36, 12).Value = 20049  # L36: 0 -> 20049
$ws.Cells.Item(36, 13).Value = -18307  # M36: -13299.667 -> -18307
$ws.Cells.Item(36, 14).Value = -21479  # N36: None -> -21479
$ws.Cells.Item(70, 8).Value = 888  # H70: 1463.6666 -> 888
$ws.Cells.Item(70, 10).Value = 0  # J70: 1751.5 -> 0
$ws.Cells.Item(70, 12).Value = 0  # L70: 5254.5 -> 0
$ws.Cells.Item(70, 14).ClearContents()  # N70: -5794.5 -> (removed)
$ws.Cells.Item(73, 8).Value = 888  # H73: 1463.6666 -> 888
$ws.Cells.Item(73, 10).Value = 0  # J73: 1751.5 -> 0
$ws.Cells.Item(73, 12).Value = 0  # L73: 5254.5 -> 0
$ws.Cells.Item(73, 14).ClearContents()  # N73: -7126.5 -> (removed)
$ws.Cells.Item(98, 8).Value = 1475  # H98: 1566.6666 -> 1475
$ws.Cells.Item(98, 9).Value = 1475  # I98: 1566.6666 -> 1475
$ws.Cells.Item(98, 11).Value = 1475  # K98: 1566.6666 -> 1475
$ws.Cells.Item(98, 13).Value = 23  # M98: -68.66660000000002 -> 23
$ws.Cells.Item(122, 8).Value = 1475  # H122: 1566.6666 -> 1475
$ws.Cells.Item(122, 9).Value = 1475  # I122: 1566.6666 -> 1475
$ws.Cells.Item(122, 11).Value = 4425  # K122: 4699.9998 -> 4425
$ws.Cells.Item(122, 13).Value = -1975  # M122: -2249.9998 -> -1975
$ws.Cells.Item(125, 8).Value = 377  # H125: 433.33334 -> 377
$ws.Cells.Item(125, 10).Value = 396.25  # J125: 500 -> 396.25
$ws.Cells.Item(125, 12).Value = 3566.25  # L125: 4500 -> 3566.25
$ws.Cells.Item(125, 14).Value = -8486.25  # N125: -9420 -> -8486.25
$ws.Cells.Item(132, 8).Value = 1222.4117  # H132: 1185.1875 -> 1222.4117
$ws.Cells.Item(132, 9).Value = 869.0769  # I132: 790 -> 869.0769
$ws.Cells.Item(132, 11).Value = 2607.2307  # K132: 2370 -> 2607.2307
$ws.Cells.Item(132, 13).Value = -77.23070000000007  # M132: 160 -> -77.23070000000007
$ws.Cells.Item(138, 8).Value = 1618.2  # H138: 534 -> 1618.2
$ws.Cells.Item(138, 9).Value = 597.5714  # I138: 534 -> 597.5714
$ws.Cells.Item(138, 10).Value = 3999.6667  # J138: 0 -> 3999.6667
$ws.Cells.Item(138, 11).Value = 1792.7142  # K138: 1602 -> 1792.7142
$ws.Cells.Item(138, 12).Value = 11999.0001  # L138: 0 -> 11999.0001
$ws.Cells.Item(138, 13).Value = 3347.2858  # M138: 3538 -> 3347.2858
$ws.Cells.Item(138, 14).Value = -22279.0001  # N138: None -> -22279.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(22, 8).Value = 14000  # H22: 0 -> 14000
$ws.Cells.Item(22, 10).Value = 14000  # J22: 0 -> 14000
$ws.Cells.Item(22, 12).Value = 14000  # L22: 0 -> 14000
$ws.Cells.Item(22, 14).Value = -14598  # N22: None -> -14598
$ws.Cells.Item(63, 8).Value = 9750.75  # H63: 11501 -> 9750.75
$ws.Cells.Item(63, 9).Value = 6332.5  # I63: 7248.75 -> 6332.5
$ws.Cells.Item(63, 11).Value = 6332.5  # K63: 7248.75 -> 6332.5
$ws.Cells.Item(63, 13).Value = -5646.5  # M63: -6562.75 -> -5646.5
$ws.Cells.Item(66, 8).Value = 9750.75  # H66: 11501 -> 9750.75
$ws.Cells.Item(66, 9).Value = 6332.5  # I66: 7248.75 -> 6332.5
$ws.Cells.Item(66, 11).Value = 31662.5  # K66: 36243.75 -> 31662.5
$ws.Cells.Item(66, 13).Value = -28230.5  # M66: -32811.75 -> -28230.5
$ws.Cells.Item(132, 8).Value = 1884.4546  # H132: 2030.1613 -> 1884.4546
$ws.Cells.Item(132, 9).Value = 933.5909  # I132: 1025.5714 -> 933.5909
$ws.Cells.Item(132, 10).Value = 3786.182  # J132: 4139.8 -> 3786.182
$ws.Cells.Item(132, 11).Value = 2800.7727  # K132: 3076.7142 -> 2800.7727
$ws.Cells.Item(132, 12).Value = 11358.546  # L132: 12419.4 -> 11358.546
$ws.Cells.Item(132, 13).Value = -270.7727  # M132: -546.7142000000003 -> -270.7727
$ws.Cells.Item(132, 14).Value = -16418.546  # N132: -17479.4 -> -16418.546

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3031.8  # H86: 3290 -> 3031.8
$ws.Cells.Item(86, 9).Value = 2666.3333  # I86: 2799.8 -> 2666.3333
$ws.Cells.Item(86, 10).Value = 3580  # J86: 4107 -> 3580
$ws.Cells.Item(86, 11).Value = 2666.3333  # K86: 2799.8 -> 2666.3333
$ws.Cells.Item(86, 12).Value = 3580  # L86: 4107 -> 3580
$ws.Cells.Item(86, 13).Value = -1543.3333  # M86: -1676.8 -> -1543.3333
$ws.Cells.Item(86, 14).Value = -5826  # N86: -6353 -> -5826
$ws.Cells.Item(89, 8).Value = 3031.8  # H89: 3290 -> 3031.8
$ws.Cells.Item(89, 9).Value = 2666.3333  # I89: 2799.8 -> 2666.3333
$ws.Cells.Item(89, 10).Value = 3580  # J89: 4107 -> 3580
$ws.Cells.Item(89, 11).Value = 13331.6665  # K89: 13999 -> 13331.6665
$ws.Cells.Item(89, 12).Value = 17900  # L89: 20535 -> 17900
$ws.Cells.Item(89, 13).Value = -7715.666499999999  # M89: -8383 -> -7715.666499999999
$ws.Cells.Item(89, 14).Value = -29132  # N89: -31767 -> -29132
$ws.Cells.Item(134, 8).Value = 2289.3333  # H134: 2870 -> 2289.3333
$ws.Cells.Item(134, 9).Value = 1267.4166  # I134: 1450 -> 1267.4166
$ws.Cells.Item(134, 10).Value = 4333.1665  # J134: 5000 -> 4333.1665
$ws.Cells.Item(134, 11).Value = 3802.2498  # K134: 4350 -> 3802.2498
$ws.Cells.Item(134, 12).Value = 12999.4995  # L134: 15000 -> 12999.4995
$ws.Cells.Item(134, 13).Value = -1267.2498  # M134: -1815 -> -1267.2498
$ws.Cells.Item(134, 14).Value = -18069.4995  # N134: -20070 -> -18069.4995

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3186.9583  # H31: 3036.077 -> 3186.9583
$ws.Cells.Item(31, 9).Value = 2265.8667  # I31: 2202.4375 -> 2265.8667
$ws.Cells.Item(31, 10).Value = 4722.1113  # J31: 4369.9 -> 4722.1113
$ws.Cells.Item(31, 11).Value = 2265.8667  # K31: 2202.4375 -> 2265.8667
$ws.Cells.Item(31, 12).Value = 4722.1113  # L31: 4369.9 -> 4722.1113
$ws.Cells.Item(31, 13).Value = -1970.8667  # M31: -1907.4375 -> -1970.8667
$ws.Cells.Item(31, 14).Value = -5312.1113  # N31: -4959.9 -> -5312.1113
$ws.Cells.Item(34, 8).Value = 3186.9583  # H34: 3036.077 -> 3186.9583
$ws.Cells.Item(34, 9).Value = 2265.8667  # I34: 2202.4375 -> 2265.8667
$ws.Cells.Item(34, 10).Value = 4722.1113  # J34: 4369.9 -> 4722.1113
$ws.Cells.Item(34, 11).Value = 2265.8667  # K34: 2202.4375 -> 2265.8667
$ws.Cells.Item(34, 12).Value = 4722.1113  # L34: 4369.9 -> 4722.1113
$ws.Cells.Item(34, 13).Value = -2063.8667  # M34: -2000.4375 -> -2063.8667
$ws.Cells.Item(34, 14).Value = -5126.1113  # N34: -4773.9 -> -5126.1113
$ws.Cells.Item(58, 8).Value = 1767.2858  # H58: 1889.8334 -> 1767.2858
$ws.Cells.Item(58, 9).Value = 1745.1666  # I58: 1887.8 -> 1745.1666
$ws.Cells.Item(58, 11).Value = 1745.1666  # K58: 1887.8 -> 1745.1666
$ws.Cells.Item(58, 13).Value = -1542.1666  # M58: -1684.8 -> -1542.1666
$ws.Cells.Item(136, 8).Value = 1767.2858  # H136: 1889.8334 -> 1767.2858
$ws.Cells.Item(136, 9).Value = 1745.1666  # I136: 1887.8 -> 1745.1666
$ws.Cells.Item(136, 11).Value = 5235.4998  # K136: 5663.4 -> 5235.4998
$ws.Cells.Item(136, 13).Value = -2685.4998  # M136: -3113.4 -> -2685.4998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(13, 8).Value = 100  # H13: 0 -> 100
$ws.Cells.Item(13, 10).Value = 100  # J13: 0 -> 100
$ws.Cells.Item(13, 12).Value = 300  # L13: 0 -> 300
$ws.Cells.Item(13, 14).Value = -636  # N13: None -> -636
$ws.Cells.Item(34, 8).Value = 3421.6428  # H34: 2616 -> 3421.6428
$ws.Cells.Item(34, 9).Value = 2999  # I34: 1749.5 -> 2999
$ws.Cells.Item(34, 10).Value = 3454.1538  # J34: 2789.3 -> 3454.1538
$ws.Cells.Item(34, 11).Value = 8997  # K34: 5248.5 -> 8997
$ws.Cells.Item(34, 12).Value = 10362.4614  # L34: 8367.900000000001 -> 10362.4614
$ws.Cells.Item(34, 13).Value = -8913  # M34: -5164.5 -> -8913
$ws.Cells.Item(34, 14).Value = -10530.4614  # N34: -8535.900000000001 -> -10530.4614
$ws.Cells.Item(39, 8).Value = 0  # H39: 18000 -> 0
$ws.Cells.Item(39, 10).Value = 0  # J39: 18000 -> 0
$ws.Cells.Item(39, 12).Value = 0  # L39: 54000 -> 0
$ws.Cells.Item(39, 14).ClearContents()  # N39: -54588 -> (removed)
$ws.Cells.Item(60, 8).Value = 2167.5  # H60: 335 -> 2167.5
$ws.Cells.Item(60, 10).Value = 4000  # J60: 0 -> 4000
$ws.Cells.Item(60, 12).Value = 12000  # L60: 0 -> 12000
$ws.Cells.Item(60, 14).Value = -12502  # N60: None -> -12502
$ws.Cells.Item(86, 8).Value = 242.5  # H86: 235.95238 -> 242.5
$ws.Cells.Item(86, 9).Value = 200  # I86: 216.66667 -> 200
$ws.Cells.Item(86, 10).Value = 270.83334  # J86: 239.16667 -> 270.83334
$ws.Cells.Item(86, 11).Value = 600  # K86: 650.00001 -> 600
$ws.Cells.Item(86, 12).Value = 812.5000200000001  # L86: 717.50001 -> 812.5000200000001
$ws.Cells.Item(86, 13).Value = 586  # M86: 535.99999 -> 586
$ws.Cells.Item(86, 14).Value = -3184.50002  # N86: -3089.50001 -> -3184.50002
$ws.Cells.Item(89, 8).Value = 242.5  # H89: 235.95238 -> 242.5
$ws.Cells.Item(89, 9).Value = 200  # I89: 216.66667 -> 200
$ws.Cells.Item(89, 10).Value = 270.83334  # J89: 239.16667 -> 270.83334
$ws.Cells.Item(89, 11).Value = 1800  # K89: 1950.00003 -> 1800
$ws.Cells.Item(89, 12).Value = 2437.50006  # L89: 2152.50003 -> 2437.50006
$ws.Cells.Item(89, 13).Value = 4128  # M89: 3977.99997 -> 4128
$ws.Cells.Item(89, 14).Value = -14293.50006  # N89: -14008.50003 -> -14293.50006
$ws.Cells.Item(109, 8).Value = 4250.92  # H109: 4403.2173 -> 4250.92
$ws.Cells.Item(109, 9).Value = 1254.6  # I109: 424.66666 -> 1254.6
$ws.Cells.Item(109, 11).Value = 3763.8  # K109: 1273.99998 -> 3763.8
$ws.Cells.Item(109, 13).Value = -2723.8  # M109: -233.9999800000001 -> -2723.8
$ws.Cells.Item(134, 8).Value = 2259.6667  # H134: 3990 -> 2259.6667
$ws.Cells.Item(134, 9).Value = 2259.6667  # I134: 3990 -> 2259.6667
$ws.Cells.Item(134, 11).Value = 6779.000100000001  # K134: 11970 -> 6779.000100000001
$ws.Cells.Item(134, 13).Value = -1709.000100000001  # M134: -6900 -> -1709.000100000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(6, 8).Value = 7669.6665  # H6: 0 -> 7669.6665
$ws.Cells.Item(6, 10).Value = 7669.6665  # J6: 0 -> 7669.6665
$ws.Cells.Item(6, 12).Value = 7669.6665  # L6: 0 -> 7669.6665
$ws.Cells.Item(6, 14).Value = -7895.6665  # N6: None -> -7895.6665
$ws.Cells.Item(13, 8).Value = 600  # H13: 1006 -> 600
$ws.Cells.Item(13, 10).Value = 600  # J13: 1006 -> 600
$ws.Cells.Item(13, 12).Value = 600  # L13: 1006 -> 600
$ws.Cells.Item(13, 14).Value = -878  # N13: -1284 -> -878
$ws.Cells.Item(16, 8).Value = 7669.6665  # H16: 0 -> 7669.6665
$ws.Cells.Item(16, 10).Value = 7669.6665  # J16: 0 -> 7669.6665
$ws.Cells.Item(16, 12).Value = 7669.6665  # L16: 0 -> 7669.6665
$ws.Cells.Item(16, 14).Value = -8169.6665  # N16: None -> -8169.6665
$ws.Cells.Item(24, 8).Value = 5755257  # H24: 23000006 -> 5755257
$ws.Cells.Item(24, 10).Value = 7007  # J24: 0 -> 7007
$ws.Cells.Item(24, 12).Value = 7007  # L24: 0 -> 7007
$ws.Cells.Item(24, 14).Value = -7353  # N24: None -> -7353
$ws.Cells.Item(80, 8).Value = 2833.3333  # H80: 3000 -> 2833.3333
$ws.Cells.Item(80, 9).Value = 2750  # I80: 3000 -> 2750
$ws.Cells.Item(80, 11).Value = 2750  # K80: 3000 -> 2750
$ws.Cells.Item(80, 13).Value = -1752  # M80: -2002 -> -1752
$ws.Cells.Item(83, 8).Value = 2833.3333  # H83: 3000 -> 2833.3333
$ws.Cells.Item(83, 9).Value = 2750  # I83: 3000 -> 2750
$ws.Cells.Item(83, 11).Value = 13750  # K83: 15000 -> 13750
$ws.Cells.Item(83, 13).Value = -8758  # M83: -10008 -> -8758
$ws.Cells.Item(122, 8).Value = 3279.3333  # H122: 2111.818 -> 3279.3333
$ws.Cells.Item(122, 9).Value = 3557.6  # I122: 2204.6667 -> 3557.6
$ws.Cells.Item(122, 10).Value = 1888  # J122: 1694 -> 1888
$ws.Cells.Item(122, 11).Value = 10672.8  # K122: 6614.000100000001 -> 10672.8
$ws.Cells.Item(122, 12).Value = 5664  # L122: 5082 -> 5664
$ws.Cells.Item(122, 13).Value = -8222.799999999999  # M122: -4164.000100000001 -> -8222.799999999999
$ws.Cells.Item(122, 14).Value = -10564  # N122: -9982 -> -10564

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 895.1667  # H16: 891.6667 -> 895.1667
$ws.Cells.Item(16, 9).Value = 833.8  # I16: 829.6 -> 833.8
$ws.Cells.Item(16, 11).Value = 833.8  # K16: 829.6 -> 833.8
$ws.Cells.Item(16, 13).Value = -663.8  # M16: -659.6 -> -663.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(31, 8).Value = 30019  # H31: 0 -> 30019
$ws.Cells.Item(31, 10).Value = 30019  # J31: 0 -> 30019
$ws.Cells.Item(31, 12).Value = 30019  # L31: 0 -> 30019
$ws.Cells.Item(31, 14).Value = -30715  # N31: None -> -30715
$ws.Cells.Item(54, 8).Value = 19999.625  # H54: 15999.9 -> 19999.625
$ws.Cells.Item(54, 10).Value = 29999.25  # J54: 29999.666 -> 29999.25
$ws.Cells.Item(54, 12).Value = 29999.25  # L54: 29999.666 -> 29999.25
$ws.Cells.Item(54, 14).Value = -31039.25  # N54: -31039.666 -> -31039.25
$ws.Cells.Item(55, 8).Value = 22000  # H55: 25000 -> 22000
$ws.Cells.Item(55, 9).Value = 20000  # I55: 0 -> 20000
$ws.Cells.Item(55, 10).Value = 24000  # J55: 25000 -> 24000
$ws.Cells.Item(55, 11).Value = 20000  # K55: 0 -> 20000
$ws.Cells.Item(55, 12).Value = 24000  # L55: 25000 -> 24000
$ws.Cells.Item(55, 13).Value = -19723  # M55: None -> -19723
$ws.Cells.Item(55, 14).Value = -24554  # N55: -25554 -> -24554
$ws.Cells.Item(61, 8).Value = 30057  # H61: 10000 -> 30057
$ws.Cells.Item(61, 9).Value = 0  # I61: 10000 -> 0
$ws.Cells.Item(61, 10).Value = 30057  # J61: 0 -> 30057
$ws.Cells.Item(61, 11).Value = 0  # K61: 10000 -> 0
$ws.Cells.Item(61, 12).Value = 30057  # L61: 0 -> 30057
$ws.Cells.Item(61, 13).ClearContents()  # M61: -9708 -> (removed)
$ws.Cells.Item(61, 14).Value = -30641  # N61: None -> -30641
$ws.Cells.Item(132, 8).Value = 1863.625  # H132: 1931.6086 -> 1863.625
$ws.Cells.Item(132, 9).Value = 915.6667  # I132: 959.6429000000001 -> 915.6667
$ws.Cells.Item(132, 11).Value = 2747.0001  # K132: 2878.9287 -> 2747.0001
$ws.Cells.Item(132, 13).Value = -217.0001000000002  # M132: -348.9287000000004 -> -217.0001000000002
